$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$cell = $ws.Cells.Item(2, 4)
$origStyle = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = '70.262.26'
$cell.Style = $origStyle

$cell = $ws.Cells.Item(2, 5)
$origStyle = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = '  +5.41%  '
$cell.Style = $origStyle

$cell = $ws.Cells.Item(3, 4)
$origStyle = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = '3.609.09'
$cell.Style = $origStyle

$cell = $ws.Cells.Item(3, 5)
$origStyle = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = '  +4.98%  '
$cell.Style = $origStyle

$cell = $ws.Cells.Item(4, 5)
$origStyle = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = '  -0.05%  '
$cell.Style = $origStyle

$cell = $ws.Cells.Item(5, 4)
$origStyle = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = '591.80'
$cell.Style = $origStyle

$cell = $ws.Cells.Item(5, 5)
$origStyle = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = '  +3.17%  '
$cell.Style = $origStyle

$cell = $ws.Cells.Item(6, 4)
$origStyle = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = '190.64'
$cell.Style = $origStyle

$cell = $ws.Cells.Item(6, 5)
$origStyle = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = '  +3.56%  '
$cell.Style = $origStyle

$cell = $ws.Cells.Item(7, 5)
$origStyle = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = '  +2.00%  '
$cell.Style = $origStyle

$cell = $ws.Cells.Item(8, 4)
$origStyle = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = '3.603.09'
$cell.Style = $origStyle

$cell = $ws.Cells.Item(8, 5)
$origStyle = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = '  +5.01%  '
$cell.Style = $origStyle

$cell = $ws.Cells.Item(9, 5)
$origStyle = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = '  -0.03%  '
$cell.Style = $origStyle

$cell = $ws.Cells.Item(10, 5)
$origStyle = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = '  +3.14%  '
$cell.Style = $origStyle

$cell = $ws.Cells.Item(11, 5)
$origStyle = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = '  +3.13%  '
$cell.Style = $origStyle

$cell = $ws.Cells.Item(12, 4)
$origStyle = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = '58.60'
$cell.Style = $origStyle

$cell = $ws.Cells.Item(12, 5)
$origStyle = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = '  +4.77%  '
$cell.Style = $origStyle

$cell = $ws.Cells.Item(13, 5)
$origStyle = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = '  +4.08%  '
$cell.Style = $origStyle

$cell = $ws.Cells.Item(14, 4)
$origStyle = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = '9.89'
$cell.Style = $origStyle

$cell = $ws.Cells.Item(15, 4)
$origStyle = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = '4.188.14'
$cell.Style = $origStyle

$cell = $ws.Cells.Item(15, 5)
$origStyle = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = '  +5.03%  '
$cell.Style = $origStyle

$cell = $ws.Cells.Item(16, 5)
$origStyle = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = '  +5.97%  '
$cell.Style = $origStyle

$cell = $ws.Cells.Item(17, 4)
$origStyle = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = '3.610.72'
$cell.Style = $origStyle

$cell = $ws.Cells.Item(17, 5)
$origStyle = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = '  +5.00%  '
$cell.Style = $origStyle

$cell = $ws.Cells.Item(18, 4)
$origStyle = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = '70.246.81'
$cell.Style = $origStyle

$cell = $ws.Cells.Item(18, 5)
$origStyle = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = '  +5.17%  '
$cell.Style = $origStyle

$cell = $ws.Cells.Item(19, 4)
$origStyle = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = '12.60'
$cell.Style = $origStyle

$cell = $ws.Cells.Item(19, 5)
$origStyle = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = '  +4.45%  '
$cell.Style = $origStyle

$cell = $ws.Cells.Item(21, 5)
$origStyle = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = '  +4.13%  '
$cell.Style = $origStyle

$cell = $ws.Cells.Item(22, 4)
$origStyle = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = '492.87'
$cell.Style = $origStyle

$cell = $ws.Cells.Item(22, 5)
$origStyle = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = '  +0.20%  '
$cell.Style = $origStyle

$cell = $ws.Cells.Item(23, 4)
$origStyle = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = '19.14'
$cell.Style = $origStyle

$cell = $ws.Cells.Item(23, 5)
$origStyle = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = '  +16.34%  '
$cell.Style = $origStyle

$cell = $ws.Cells.Item(24, 4)
$origStyle = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = '5.36'
$cell.Style = $origStyle

$cell = $ws.Cells.Item(24, 5)
$origStyle = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = '  +6.64%  '
$cell.Style = $origStyle

$cell = $ws.Cells.Item(25, 4)
$origStyle = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = '4.47'
$cell.Style = $origStyle

$cell = $ws.Cells.Item(25, 5)
$origStyle = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = '  +5.16%  '
$cell.Style = $origStyle

$cell = $ws.Cells.Item(26, 4)
$origStyle = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = '90.96'
$cell.Style = $origStyle

$cell = $ws.Cells.Item(26, 5)
$origStyle = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = '  +1.38%  '
$cell.Style = $origStyle

$cell = $ws.Cells.Item(27, 5)
$origStyle = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = '  +6.11%  '
$cell.Style = $origStyle

$cell = $ws.Cells.Item(28, 4)
$origStyle = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = '11.19'
$cell.Style = $origStyle

$cell = $ws.Cells.Item(28, 5)
$origStyle = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = '  +1.13%  '
$cell.Style = $origStyle

$cell = $ws.Cells.Item(29, 4)
$origStyle = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = '9.60'
$cell.Style = $origStyle

$cell = $ws.Cells.Item(29, 5)
$origStyle = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = '  +4.87%  '
$cell.Style = $origStyle

$cell = $ws.Cells.Item(30, 4)
$origStyle = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = '32.93'
$cell.Style = $origStyle

$cell = $ws.Cells.Item(30, 5)
$origStyle = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = '  +4.94%  '
$cell.Style = $origStyle

$cell = $ws.Cells.Item(31, 4)
$origStyle = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = '7.74'
$cell.Style = $origStyle

$cell = $ws.Cells.Item(31, 5)
$origStyle = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = '  +8.62%  '
$cell.Style = $origStyle

$cell = $ws.Cells.Item(32, 4)
$origStyle = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = '635.48'
$cell.Style = $origStyle

$cell = $ws.Cells.Item(32, 5)
$origStyle = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = '  +6.99%  '
$cell.Style = $origStyle

$cell = $ws.Cells.Item(33, 4)
$origStyle = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = '12.32'
$cell.Style = $origStyle

$cell = $ws.Cells.Item(33, 5)
$origStyle = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = '  +5.51%  '
$cell.Style = $origStyle

$cell = $ws.Cells.Item(34, 5)
$origStyle = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = '  +6.89%  '
$cell.Style = $origStyle

$cell = $ws.Cells.Item(35, 4)
$origStyle = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = '65.79'
$cell.Style = $origStyle

$cell = $ws.Cells.Item(35, 5)
$origStyle = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = '  +4.15%  '
$cell.Style = $origStyle

$cell = $ws.Cells.Item(36, 4)
$origStyle = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = '38.71'
$cell.Style = $origStyle

$cell = $ws.Cells.Item(36, 5)
$origStyle = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = '  +6.82%  '
$cell.Style = $origStyle

$cell = $ws.Cells.Item(37, 4)
$origStyle = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = '0.0₃0819'
$cell.Style = $origStyle

$cell = $ws.Cells.Item(37, 5)
$origStyle = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = '  +6.09%  '
$cell.Style = $origStyle

$cell = $ws.Cells.Item(38, 4)
$origStyle = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = '0.406'
$cell.Style = $origStyle

$cell = $ws.Cells.Item(38, 5)
$origStyle = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = '  +5.50%  '
$cell.Style = $origStyle

$cell = $ws.Cells.Item(39, 5)
$origStyle = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = '  +0.02%  '
$cell.Style = $origStyle

$cell = $ws.Cells.Item(40, 5)
$origStyle = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = '  -0.71%  '
$cell.Style = $origStyle

$cell = $ws.Cells.Item(41, 4)
$origStyle = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = '3.56'
$cell.Style = $origStyle

$cell = $ws.Cells.Item(42, 4)
$origStyle = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = '3.301.63'
$cell.Style = $origStyle

$cell = $ws.Cells.Item(42, 5)
$origStyle = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = '  +4.06%  '
$cell.Style = $origStyle

$cell = $ws.Cells.Item(43, 4)
$origStyle = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = '3.15'
$cell.Style = $origStyle

$cell = $ws.Cells.Item(43, 5)
$origStyle = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = '  +7.57%  '
$cell.Style = $origStyle

$cell = $ws.Cells.Item(44, 2)
$origStyle = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = 'Fetch.AI'
$cell.Style = $origStyle

$cell = $ws.Cells.Item(44, 3)
$origStyle = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = 'https://coinranking.com/coin/AWma-WzFHmKVQ+fetchai-fet'
$cell.Style = $origStyle

$cell = $ws.Cells.Item(44, 4)
$origStyle = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = '2.75'
$cell.Style = $origStyle

$cell = $ws.Cells.Item(44, 5)
$origStyle = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = '  +8.18%  '
$cell.Style = $origStyle

$cell = $ws.Cells.Item(45, 2)
$origStyle = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = 'VeChain'
$cell.Style = $origStyle

$cell = $ws.Cells.Item(45, 3)
$origStyle = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = 'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet'
$cell.Style = $origStyle

$cell = $ws.Cells.Item(45, 4)
$origStyle = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = '0.0453'
$cell.Style = $origStyle

$cell = $ws.Cells.Item(45, 5)
$origStyle = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = '  +5.54%  '
$cell.Style = $origStyle

$cell = $ws.Cells.Item(46, 5)
$origStyle = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = '  +2.86%  '
$cell.Style = $origStyle

$cell = $ws.Cells.Item(47, 4)
$origStyle = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = '3.27'
$cell.Style = $origStyle

$cell = $ws.Cells.Item(47, 5)
$origStyle = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = '  +2.10%  '
$cell.Style = $origStyle

$cell = $ws.Cells.Item(48, 4)
$origStyle = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = '9.07'
$cell.Style = $origStyle

$cell = $ws.Cells.Item(48, 5)
$origStyle = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = '  +3.50%  '
$cell.Style = $origStyle

$cell = $ws.Cells.Item(49, 5)
$origStyle = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = '  -2.57%  '
$cell.Style = $origStyle

$cell = $ws.Cells.Item(50, 4)
$origStyle = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = '3.31'
$cell.Style = $origStyle

$cell = $ws.Cells.Item(50, 5)
$origStyle = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = '  +5.00%  '
$cell.Style = $origStyle

$cell = $ws.Cells.Item(51, 2)
$origStyle = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = 'Monero'
$cell.Style = $origStyle

$cell = $ws.Cells.Item(51, 3)
$origStyle = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = 'https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr'
$cell.Style = $origStyle

$cell = $ws.Cells.Item(51, 4)
$origStyle = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = '143.14'
$cell.Style = $origStyle

$cell = $ws.Cells.Item(51, 5)
$origStyle = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = '  +2.04%  '
$cell.Style = $origStyle